# Generate Report for Handback
# Updates the localization-status workbook so that the zh-cn/de-de "Status"
# columns show the handed-back state, fills in the "Latest Target File" /
# "Latest Handback File" (and, where the handback has completed, the
# "Latest Handback DateTime") columns with links/info for each localized
# file, and widens the columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$c9Name = "c9f865dc-995e-4b48-bea7-cc95fc68ab88.md"
$caName = "ca6b369d-cf96-4d59-bb01-34676d94e389.md"
$c9Url  = "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/324f98ad6dc7c0ab3e82e00265beaf5be8c6bada/e2e/c9f865dc-995e-4b48-bea7-cc95fc68ab88.md"
$caUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/324f98ad6dc7c0ab3e82e00265beaf5be8c6bada/e2e/ca6b369d-cf96-4d59-bb01-34676d94e389.md"

$c9ZhXlf = "c9f865dc-995e-4b48-bea7-cc95fc68ab88.b8e53fba51eae758d8bffaba39ad430444cb9fa6.zh-cn.xlf"
$caZhXlf = "ca6b369d-cf96-4d59-bb01-34676d94e389.8a5cae30e2caaf7ae8d41157ebe3f963f6313043.zh-cn.xlf"
$c9DeXlf = "c9f865dc-995e-4b48-bea7-cc95fc68ab88.b8e53fba51eae758d8bffaba39ad430444cb9fa6.de-de.xlf"
$caDeXlf = "ca6b369d-cf96-4d59-bb01-34676d94e389.8a5cae30e2caaf7ae8d41157ebe3f963f6313043.de-de.xlf"

$zhHandbackDateTime = "2017-01-03 06:23:30"
$deHandbackDateTime = "2017-01-03 06:23:41"

# ---- Overview sheet: update status + widen the per-language columns ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus
$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(10).ColumnWidth = 40
$wsZh.Columns.Item(11).ColumnWidth = 40

$wsZh.Range("J2").Value = $c9Name
$wsZh.Hyperlinks.Add($wsZh.Range("J2"), $c9Url, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $c9Name)
$wsZh.Range("K2").Value = $c9ZhXlf
$wsZh.Range("L2").Value = $zhHandbackDateTime

$wsZh.Range("J3").Value = $caName
$wsZh.Hyperlinks.Add($wsZh.Range("J3"), $caUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $caName)
$wsZh.Range("K3").Value = $caZhXlf
$wsZh.Range("L3").Value = $zhHandbackDateTime

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus
$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(10).ColumnWidth = 40
$wsDe.Columns.Item(11).ColumnWidth = 40

$wsDe.Range("J2").Value = $c9Name
$wsDe.Hyperlinks.Add($wsDe.Range("J2"), $c9Url, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $c9Name)
$wsDe.Range("K2").Value = $c9DeXlf
$wsDe.Range("L2").Value = $deHandbackDateTime

$wsDe.Range("J3").Value = $caName
$wsDe.Hyperlinks.Add($wsDe.Range("J3"), $caUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $caName)
$wsDe.Range("K3").Value = $caDeXlf
$wsDe.Range("L3").Value = $deHandbackDateTime
